$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 2 ("In recent years, ...") - split the final run so that "Gov"
# and the trailing "issues." get wrapped in <w:proofErr> spell/grammar marks,
# matching Word's automatic proofing-error bookmarks.
# ---------------------------------------------------------------------------
$para2Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00FC5B51" w:rsidRDefault="00FC5B51" w:rsidP="00FC5B51"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>In recent years, driverless cars have moved from an idea firmly planted in science fiction into something much more real and achievable by modern technology. With Google&#8217;s self-driving cars having travelled 700,000 miles as of August 2014 (Gomes,</w:t></w:r><w:r w:rsidR="00360885"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>2014) and fully driverless cars anticipated to be brought onto the UK market within the 2020s (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Gov</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">, 2015), driverless cars could be set to take over the motor industry. However, with any such advanced technology there are a significant number of ethical issues which will potentially hinder the progress of these vehicles making it to market and this paper aims to explore such </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>issues.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
'@

$d.Paragraphs(2).Range.InsertXML($para2Xml)

# ---------------------------------------------------------------------------
# Paragraph 4 (the Trolley Problem paragraph) - split the first run so that
# ",1985" is wrapped in <w:proofErr w:type="gramStart/gramEnd"/>, append the
# "[PLEASE WRITE SOMETHING HERE ...]" note, and drop the _GoBack bookmark
# (it moves into its own paragraph below).
# ---------------------------------------------------------------------------
$para4Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="004B13FD" w:rsidRPr="00FC5B51" w:rsidRDefault="004B13FD" w:rsidP="00FC5B51"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>A well-known thought experiment which covers these issues is the Trolley Problem (Thomson</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>,1985</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>). The problem is described as follows: Imagine you are driving a trolley and you come to a point where the track forks in two directions. In one direction you find 5 people, unable to move out of the way in time, and on the other side only 1 person in the same position. You find your brakes are broken so you are guaranteed to kill at least 1 person. Is it the correct decision to turn the trolley onto the side with 1 person? This problem is clearly applicable to the aforementioned programming of these cars. Most people would take the utilitarian approach in this scenario; opt to kill the lone person as it maximises the total lives saved. An argument can be made against such an approach though, as has been discussed by Ian Chapman (2</w:t></w:r><w:r w:rsidR="00360885"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">015). He brings up the issue of a slightly more subtle application of the trolley problem. Suppose a driverless car has a decision between hitting (and potentially injuring) one of two motorcyclists on the road: </w:t></w:r><w:r w:rsidR="00360885"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:lastRenderedPageBreak/><w:t>the first motorcyclist is wearing a helmet whereas the second is not. If the car were to act according to utilitarianism then it would hit the cyclist wearing the helmet, which could seem questionable to some as, if taken as a general rule, would suggest that people attempting to put themselves out of harm&#8217;s way, are in fact put directly at harm in such conditions.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> [PLEASE WRITE SOMETHING HERE REGARDING HOW THIS COULD BE SENSIBLY DEALT WITH &#8211; MORE RESEARCH COULD BE REQUIRED]</w:t></w:r></w:p>
'@

$d.Paragraphs(4).Range.InsertXML($para4Xml)

# ---------------------------------------------------------------------------
# Add the two new paragraphs at the end of the document: the new "If
# driverless cars..." discussion, then an (otherwise empty) paragraph that
# now owns the _GoBack bookmark.
# ---------------------------------------------------------------------------
$d.Paragraphs(4).Range.InsertParagraphAfter()
$d.Paragraphs(5).Range.InsertParagraphAfter()

$para5Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">If driverless cars were to become fully autonomous then one can view them as an independent decision maker. If we assume that the car doesn&#8217;t treat any </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/></w:rPr><w:t xml:space="preserve">human </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>life more important than another (a different case may be applicable for animals etc.) then the car could potentially decide to kill the driver (or more accurately in this case, passenger) if it appeared to be the best option. At this stage, it is notable that this may not be the same decision as a human driver would make so it is worth questioning: if we are questioning the safety and decision making of driverless cars, shouldn&#8217;t we question human decision making before using it as a standard for autonomous vehicles?</w:t></w:r></w:p>
'@

$d.Paragraphs(5).Range.InsertXML($para5Xml)

$para6Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$d.Paragraphs(6).Range.InsertXML($para6Xml)

Write-Output "done"
